$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.241.39"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.598.31"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.823.97"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.611.86"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "26.253.16"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.91%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "1.444.43"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.929"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "1.736.23"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.755"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0947"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.13%  "
